# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.521.79"
$ws.Range("E2").Value = "  -7.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.209.32"
$ws.Range("E3").Value = "  -8.03%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.69"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.52"
$ws.Range("E6").Value = "  -14.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -9.96%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  -11.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.50"
$ws.Range("E10").Value = "  -12.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.55"
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0829"
$ws.Range("E12").Value = "  -11.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.55"
$ws.Range("E13").Value = "  -13.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.104"
$ws.Range("E14").Value = "  -4.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.860"
$ws.Range("E15").Value = "  -14.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.552.55"
$ws.Range("E16").Value = "  -7.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.84"
$ws.Range("E17").Value = "  -12.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.214.32"
$ws.Range("E18").Value = "  -7.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.458.24"
$ws.Range("E19").Value = "  -7.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.54"
$ws.Range("E21").Value = "  -12.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0943"
$ws.Range("E22").Value = "  -13.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.17"
$ws.Range("E23").Value = "  -10.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.53"
$ws.Range("E24").Value = "  -13.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "232.67"
$ws.Range("E25").Value = "  -12.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.10"
$ws.Range("E26").Value = "  -10.58%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -11.30%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("E29").Value = "  -8.75%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.40"
$ws.Range("E30").Value = "  -15.98%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.31"
$ws.Range("E31").Value = "  -10.97%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0869"
$ws.Range("E32").Value = "  -11.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.52"
$ws.Range("E33").Value = "  -9.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.72"
$ws.Range("E34").Value = "  -16.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.70"
$ws.Range("E35").Value = "  -8.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.20"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("E37").Value = "  -8.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.36"
$ws.Range("E38").Value = "  -11.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("E39").Value = "  -12.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0318"
$ws.Range("E41").Value = "  -12.24%  "
$ws.Range("E42").Value = "  -15.66%  "
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.764.30"
$ws.Range("E44").Value = "  +6.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "87.94"
$ws.Range("E45").Value = "  -14.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.81"
$ws.Range("E46").Value = "  -13.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.204"
$ws.Range("E47").Value = "  -15.15%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.32"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.60"
$ws.Range("E49").Value = "  -14.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.82"
$ws.Range("E50").Value = "  -16.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.47"
$ws.Range("E51").Value = "  -10.73%  "
